# Consolidate prompts: remove DISCUSSION_SUMMARY and TEACHER_FEEDBACK_EXTRACTION rows/prompts,
# drop the discussion_summary column from Discussions and the teacher_feedback column from
# Transcripts, and update the GROUP_FEEDBACK / INDIVIDUAL_FEEDBACK prompt text accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Discussions: remove column I ("discussion_summary"). Columns J:O shift left
# to I:N (approved, canvas_assignment_id, canvas_item_type, error_message,
# created_at, updated_at) along with their widths.
# ---------------------------------------------------------------------------
$wsDiscussions = $wb.Worksheets.Item("Discussions")
$wsDiscussions.Columns.Item(9).Delete()

# ---------------------------------------------------------------------------
# Transcripts: remove column E ("teacher_feedback"). Columns F:G shift left
# to E:F (created_at, updated_at) along with their widths.
# ---------------------------------------------------------------------------
$wsTranscripts = $wb.Worksheets.Item("Transcripts")
$wsTranscripts.Columns.Item(5).Delete()

# ---------------------------------------------------------------------------
# Prompts: drop the DISCUSSION_SUMMARY row (row 3) and the
# TEACHER_FEEDBACK_EXTRACTION row (originally row 6, now row 5 after the
# first deletion shifts everything up by one).
# ---------------------------------------------------------------------------
$wsPrompts = $wb.Worksheets.Item("Prompts")
$wsPrompts.Rows.Item(3).Delete()
$wsPrompts.Rows.Item(5).Delete()

# Narrow column A now that TEACHER_FEEDBACK_EXTRACTION (the longest prompt_name)
# is gone.
$wsPrompts.Columns.Item(1).ColumnWidth = 25.166666666666668

# GROUP_FEEDBACK is now row 3 - drop the "Teacher's notes: {teacher_feedback}"
# line (the prompt already receives the full transcript and is instructed to
# look for the teacher's oral feedback within it).
$groupFeedback = @'
You are a high school teacher analyzing a Harkness discussion. You will produce exactly two paragraphs.

**PARAGRAPH 1 — Discussion Summary** (Neutral Voice)
Write in a neutral, objective, third-person voice. Provide a detailed summary of the discussion's main topics and flow. Identify 2-3 "defining moments" — key turning points, breakthrough ideas, or significant challenges that shaped the conversation.

**PARAGRAPH 2 — Evaluative Comment** (Teacher Voice)
Write in the teacher's voice, directed at the class ("you" plural, "I" for the teacher). The tone must be direct, informal, supportive, and clear. Follow this mandatory "Critique Sandwich" structure:

1. **The Grade**: State the grade clearly and colloquially in the first sentence. (e.g., "This was a strong discussion, earning a solid 8.5 out of 10.", "This was a decent but not great start... 7/10.")
2. **The Good**: Highlight 2-3 specific positive achievements. Credit specific students by name, linking them to their idea or contribution.
3. **The Gap**: Identify the primary weakness or area for growth.
4. **The Next Step**: Conclude with a single, clear, actionable goal for the next discussion.

**Tone alignment with grade:**
- High grade (9-10): Frame positives as "excellent" or "deep"; the gap is a "final step" to the next level.
- Medium grade (7-8.5): Balanced ("solid," "decent start") with a more significant gap to work on.
- Lower grade (below 7): Honest but encouraging; clear gap with concrete next steps.

**Important:**
- If the teacher gave oral feedback during the discussion (often near the end — look for phrases like "my evaluation," "my feedback," or the teacher summarizing), align your evaluation with their points.
- Credit specific students by name for notable contributions.
- If the teacher intervened to guide the discussion, acknowledge this (e.g., "I had to provide the key synthesizing question").

Grade: {grade}

Transcript:
{transcript}

Write the two paragraphs now (summary paragraph first, then evaluative comment):
'@
$wsPrompts.Range("B3").Value = $groupFeedback
# Re-assigning a wrapped-text cell auto-grows the row to fit; put the
# explicit 200pt height (unchanged by this edit) back.
$wsPrompts.Rows.Item(3).RowHeight = 200

# INDIVIDUAL_FEEDBACK is now row 4 - drop the same "Teacher's notes" line and
# add an explicit "**Important:**" instruction to look for the teacher's oral
# feedback in the transcript instead.
$individualFeedback = @'
You are a high school teacher providing personalized feedback to {student_name} about their Harkness discussion participation. You will produce exactly two paragraphs.

**PARAGRAPH 1 — Contribution Summary** (Neutral Voice)
Write in a neutral, objective voice. Summarize what {student_name} contributed to the discussion — their main points, arguments, and how they engaged with other students' ideas. Note specific moments where they advanced or redirected the conversation.

**PARAGRAPH 2 — Evaluative Comment** (Teacher Voice)
Write in the teacher's voice, directed at the student ("you"). The tone must be direct, informal, supportive, and clear. Follow this "Critique Sandwich" structure:

1. **The Grade**: State the grade clearly in the first sentence.
2. **The Good**: Highlight 2-3 specific strengths from their participation, referencing actual points they made.
3. **The Gap**: Identify their primary area for growth as a discussion participant.
4. **The Next Step**: Conclude with a single, actionable goal for their next discussion.

**Tone alignment with grade:**
- High grade (9-10): "Excellent" contributions; the gap is a stretch goal.
- Medium grade (7-8.5): "Solid" participation with clear room to grow.
- Lower grade (below 7): Encouraging but honest about what's missing.

**Important:**
- If the teacher gave oral feedback during the discussion (often near the end — look for phrases like "my evaluation," "my feedback," or the teacher summarizing), align your evaluation with their points.

Grade: {grade}

{student_name}'s contributions:
{contributions}

Full discussion transcript (for context):
{transcript}

Write the two paragraphs now (contribution summary first, then evaluative comment for {student_name}):
'@
$wsPrompts.Range("B4").Value = $individualFeedback
# Same autofit correction as row 3 above.
$wsPrompts.Rows.Item(4).RowHeight = 200

Write-Output "edit complete"
